$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "ESFUERZO REAL" (actual effort) values for days 2 (F) and 3 (G)
# across every task row, and clear the now-stale "Responsable" (M column)
# value for the tasks that were finished (their work is now logged so the
# "still pending / assigned" marker is removed).

# Tarea 1 - Finalizar una jornada y guardar los datos
$ws.Range("F8").Value = 60
$ws.Range("G8").Value = 60
$ws.Range("M8").Value = ""

# Tarea 2 - Iniciar la jornada, iniciar la primera actividad
$ws.Range("F10").Value = 120
$ws.Range("G10").Value = 45
$ws.Range("M10").Value = ""

# Tarea 2.1 - Creacion de la base de datos
$ws.Range("F11").Value = 80
$ws.Range("G11").Value = 80
$ws.Range("M11").Value = ""

# Tarea 2.2 - Creacion de los botones necesarios
$ws.Range("F12").Value = 80
$ws.Range("G12").Value = 80
$ws.Range("M12").Value = ""

# Tarea 2.3 - Lógica de la aplicación
$ws.Range("F13").Value = 60
$ws.Range("G13").Value = 60
$ws.Range("M13").Value = ""

# Tarea 6 - Creación de la ayuda
$ws.Range("F16").Value = 60
$ws.Range("G16").Value = 60
$ws.Range("M16").Value = ""

# Move the active selection to match the author's final cursor position
$ws.Range("G10").Select()
